$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Events-Tillage")

# "Input [Source Name]" -> "Input [Sample Name]"
$ws.Range("A1").Value = "Input [Sample Name]"

# "Output [Source Name]" -> "Output [Sample Name]"
$ws.Range("R1").Value = "Output [Sample Name]"
